$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.004.00"
$ws.Range("E2").Value = "  -0.73%  "
$ws.Range("D3").Value = "3.539.26"
$ws.Range("E3").Value = "  -0.86%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "618.68"
$ws.Range("E5").Value = "  +6.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "185.66"
$ws.Range("E6").Value = "  +1.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.631"
$ws.Range("E7").Value = "  +1.66%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.214"
$ws.Range("E9").Value = "  -2.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.662"
$ws.Range("E10").Value = "  +2.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.79"
$ws.Range("E11").Value = "  -0.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000308"
$ws.Range("E12").Value = "  -4.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.58"
$ws.Range("E13").Value = "  +0.93%  "
$ws.Range("D14").Value = "4.105.88"
$ws.Range("E14").Value = "  -0.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "630.13"
$ws.Range("E15").Value = "  +10.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "12.86"
$ws.Range("E16").Value = "  +4.32%  "
$ws.Range("D17").Value = "70.106.72"
$ws.Range("E17").Value = "  -0.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.01"
$ws.Range("E18").Value = "  -2.48%  "
$ws.Range("D19").Value = "3.551.23"
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.997"
$ws.Range("E21").Value = "  -1.03%  "
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.75"
$ws.Range("E23").Value = "  +3.98%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "103.35"
$ws.Range("E24").Value = "  +8.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.93"
$ws.Range("E25").Value = "  -0.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.04"
$ws.Range("E26").Value = "  +3.87%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.00"
$ws.Range("E27").Value = "  -2.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.83"
$ws.Range("E28").Value = "  +8.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "34.62"
$ws.Range("E29").Value = "  +7.62%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.08"
$ws.Range("E30").Value = "  -2.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.43"
$ws.Range("E31").Value = "  +1.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.116"
$ws.Range("E32").Value = "  +1.62%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "63.89"
$ws.Range("E33").Value = "  -0.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.65"
$ws.Range("E34").Value = "  +19.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.24"
$ws.Range("E35").Value = "  -2.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "530.48"
$ws.Range("E36").Value = "  -4.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.401"
$ws.Range("E37").Value = "  -2.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.23"
$ws.Range("E39").Value = "  -0.62%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.57"
$ws.Range("E40").Value = "  +6.76%  "
$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").Value = "0.0₃0778"
$ws.Range("E41").Value = "  -4.07%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "3.527.41"
$ws.Range("E42").Value = "  +3.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.136"
$ws.Range("E43").Value = "  +0.63%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0463"
$ws.Range("E44").Value = "  +4.71%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.96"
$ws.Range("E45").Value = "  +0.40%  "
$ws.Range("E46").Value = "  +4.38%  "
$ws.Range("E47").Value = "  -5.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.08"
$ws.Range("E48").Value = "  -3.65%  "
$ws.Range("E49").Value = "  +0.44%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.41"
$ws.Range("E50").Value = "  -2.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "134.03"
$ws.Range("E51").Value = "  -1.42%  "

Write-Output "Applied cryptos update"
